# Assign DEBITAMT and CREDITAMT on Batch_Detail based on positive/negative DTLAMOUNT,
# and update BATCHID from 536924 to 825001 on both Batch_Header and Batch_Detail.

$wb = $excel.ActiveWorkbook

# --- Batch_Header: update BATCHID (column B) for data rows 2-5 ---
$wsHeader = $wb.Worksheets.Item("Batch_Header")
$headerLastRow = $wsHeader.UsedRange.Rows.Count
for ($r = 2; $r -le $headerLastRow; $r++) {
    $batchId = $wsHeader.Cells.Item($r, 2).Value2
    if ($batchId -eq 536924) {
        $wsHeader.Cells.Item($r, 2).Value2 = 825001
    }
}

# --- Batch_Detail: update BATCHID (column C) and DEBITAMT/CREDITAMT (AF/AG) ---
$wsDetail = $wb.Worksheets.Item("Batch_Detail")
$detailLastRow = $wsDetail.UsedRange.Rows.Count

# Column indices: C=3 (BATCHID), M=13 (DTLAMOUNT), AF=32 (DEBITAMT), AG=33 (CREDITAMT)
$colBatchId = 3
$colDtlAmount = 13
$colDebitAmt = 32
$colCreditAmt = 33

for ($r = 2; $r -le $detailLastRow; $r++) {
    $batchId = $wsDetail.Cells.Item($r, $colBatchId).Value2
    if ($batchId -eq 536924) {
        $wsDetail.Cells.Item($r, $colBatchId).Value2 = 825001
    }

    $dtlAmount = $wsDetail.Cells.Item($r, $colDtlAmount).Value2
    if ($dtlAmount -lt 0) {
        $wsDetail.Cells.Item($r, $colDebitAmt).Value2 = 0
        $wsDetail.Cells.Item($r, $colCreditAmt).Value2 = $dtlAmount
    } else {
        $wsDetail.Cells.Item($r, $colDebitAmt).Value2 = $dtlAmount
        $wsDetail.Cells.Item($r, $colCreditAmt).Value2 = 0
    }
}
